# Adding WY 2024 (and WY 2025) ice-off data rows to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: water year 2024 ice-off date (6/6/2024), day-of-year 250.
$ws.Range("A13").Value = 45449
$ws.Range("A13").NumberFormat = "mm-dd-yy"
$ws.Range("B13").Value = 250
$ws.Range("C13").Value = 2024

# Row 14: water year 2025 ice-off date (5/6/2025), day-of-year 218.
# Reuse the same date number format already used by the rest of column A.
$ws.Range("A14").Value = 45783
$ws.Range("A14").NumberFormat = $ws.Range("A12").NumberFormat()
$ws.Range("B14").Value = 218
$ws.Range("C14").Value = 2025

# Leave the active selection on D14, matching the saved workbook state.
$ws.Range("D14").Select()
